# This script applies the weekly data shift described by the commit.
# Rows 47-151: columns D (Fecha), I (Calidad), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), O (Origen), P (Precio $/Kg)
# each take on the values previously held by the row above them (row N-1),
# row 47 receives a brand-new weekly record, and a new row 152 is appended
# holding the values that used to belong to the (old) last row, 151.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$records = @(
    @{Row=47; D=44519; I="Primera"; J=200; K=600; L=700; M=650; O="Provincia de Diguillín"; P=650},
    @{Row=48; D=44392; I="Primera"; J=300; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=49; D=44489; I="Primera"; J=200; K=600; L=700; M=650; O="Región del Maule"; P=650},
    @{Row=50; D=44434; I="Primera"; J=300; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=51; D=44449; I="Primera"; J=300; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=52; D=44399; I="Primera"; J=120; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=53; D=44298; I="Primera"; J=120; K=800; L=850; M=825; O="Provincia de Diguillín"; P=825},
    @{Row=54; D=44482; I="Primera"; J=300; K=600; L=700; M=650; O="Provincia de Diguillín"; P=650},
    @{Row=55; D=44405; I="Primera"; J=160; K=550; L=600; M=575; O="Provincia de Diguillín"; P=575},
    @{Row=56; D=44218; I="Primera"; J=2600; K=800; L=900; M=846; O="Región del Maule"; P=846},
    @{Row=57; D=44273; I="Primera"; J=3300; K=950; L=1000; M=977; O="Región del Maule"; P=977},
    @{Row=58; D=44386; I="Primera"; J=160; K=700; L=750; M=725; O="Provincia de Diguillín"; P=725},
    @{Row=59; D=44435; I="Primera"; J=1200; K=500; L=650; M=600; O="Provincia de Diguillín"; P=600},
    @{Row=60; D=44328; I="Primera"; J=300; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=61; D=44442; I="Primera"; J=300; K=650; L=700; M=675; O="Región del Maule"; P=675},
    @{Row=62; D=44516; I="Primera"; J=300; K=600; L=700; M=650; O="Provincia de Diguillín"; P=650},
    @{Row=63; D=44175; I="Primera"; J=120; K=750; L=800; M=775; O="Provincia de Diguillín"; P=775},
    @{Row=64; D=44168; I="Primera"; J=120; K=700; L=750; M=725; O="Provincia de Diguillín"; P=725},
    @{Row=65; D=44203; I="Primera"; J=2800; K=700; L=800; M=754; O="Provincia de Diguillín"; P=754},
    @{Row=66; D=44475; I="Primera"; J=300; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=67; D=44200; I="Primera"; J=2600; K=700; L=800; M=746; O="Provincia de Diguillín"; P=746},
    @{Row=68; D=44419; I="Primera"; J=300; K=700; L=750; M=725; O="Provincia de Diguillín"; P=725},
    @{Row=69; D=44162; I="Primera"; J=2800; K=900; L=1000; M=946; O="Región de Coquimbo"; P=946},
    @{Row=70; D=44357; I="Primera"; J=300; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=71; D=44202; I="Primera"; J=3200; K=700; L=800; M=753; O="Provincia de Diguillín"; P=753},
    @{Row=72; D=44390; I="Primera"; J=120; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=73; D=44174; I="Primera"; J=120; K=800; L=850; M=825; O="Provincia de Diguillín"; P=825},
    @{Row=74; D=44293; I="Primera"; J=300; K=800; L=850; M=825; O="Provincia de Diguillín"; P=825},
    @{Row=75; D=44496; I="Primera"; J=200; K=600; L=700; M=650; O="Región del Maule"; P=650},
    @{Row=76; D=44326; I="Primera"; J=300; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=77; D=44302; I="Primera"; J=300; K=750; L=800; M=775; O="Provincia de Diguillín"; P=775},
    @{Row=78; D=44308; I="Primera"; J=300; K=700; L=750; M=725; O="Provincia de Diguillín"; P=725},
    @{Row=79; D=44498; I="Primera"; J=300; K=650; L=700; M=675; O="Provincia de Diguillín"; P=675},
    @{Row=80; D=44420; I="Primera"; J=120; K=600; L=650; M=625; O="Región del Maule"; P=625},
    @{Row=81; D=44398; I="Primera"; J=120; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=82; D=44396; I="Primera"; J=240; K=500; L=600; M=550; O="Provincia de Diguillín"; P=550},
    @{Row=83; D=44321; I="Primera"; J=300; K=700; L=750; M=725; O="Provincia de Diguillín"; P=725},
    @{Row=84; D=44208; I="Primera"; J=2900; K=700; L=800; M=748; O="Región del Maule"; P=748},
    @{Row=85; D=44349; I="Primera"; J=300; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=86; D=44477; I="Primera"; J=300; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=87; D=44487; I="Primera"; J=100; K=600; L=700; M=650; O="Provincia de Diguillín"; P=650},
    @{Row=88; D=44452; I="Primera"; J=300; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=89; D=44505; I="Primera"; J=400; K=600; L=700; M=650; O="Provincia de Diguillín"; P=650},
    @{Row=90; D=44204; I="Primera"; J=2800; K=750; L=800; M=773; O="Provincia de Diguillín"; P=773},
    @{Row=91; D=44306; I="Primera"; J=300; K=750; L=800; M=775; O="Provincia de Diguillín"; P=775},
    @{Row=92; D=44509; I="Primera"; J=400; K=600; L=700; M=650; O="Provincia de Diguillín"; P=650},
    @{Row=93; D=44454; I="Primera"; J=300; K=500; L=600; M=550; O="Provincia de Diguillín"; P=550},
    @{Row=94; D=44278; I="Primera"; J=300; K=800; L=850; M=825; O="Provincia de Diguillín"; P=825},
    @{Row=95; D=44265; I="Primera"; J=2700; K=1100; L=1200; M=1156; O="Región Metropolitana"; P=1156},
    @{Row=96; D=44494; I="Primera"; J=200; K=600; L=700; M=650; O="Región del Maule"; P=650},
    @{Row=97; D=44300; I="Primera"; J=160; K=800; L=850; M=825; O="Provincia de Diguillín"; P=825},
    @{Row=98; D=44356; I="Primera"; J=300; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=99; D=44469; I="Primera"; J=300; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=100; D=44446; I="Primera"; J=300; K=600; L=650; M=625; O="Región del Maule"; P=625},
    @{Row=101; D=44323; I="Primera"; J=400; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=102; D=44417; I="Primera"; J=300; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=103; D=44342; I="Primera"; J=300; K=700; L=750; M=725; O="Provincia de Diguillín"; P=725},
    @{Row=104; D=44406; I="Primera"; J=120; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=105; D=44295; I="Primera"; J=120; K=800; L=850; M=825; O="Provincia de Diguillín"; P=825},
    @{Row=106; D=44270; I="Primera"; J=2700; K=1000; L=1100; M=1048; O="Región del Maule"; P=1048},
    @{Row=107; D=44363; I="Primera"; J=120; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=108; D=44299; I="Primera"; J=160; K=700; L=750; M=725; O="Provincia de Diguillín"; P=725},
    @{Row=109; D=44372; I="Primera"; J=300; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=110; D=44372; I="Segunda"; J=80; K=500; L=500; M=500; O="Provincia de Diguillín"; P=500},
    @{Row=111; D=44403; I="Primera"; J=300; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=112; D=44169; I="Primera"; J=120; K=750; L=800; M=775; O="Provincia de Diguillín"; P=775},
    @{Row=113; D=44195; I="Primera"; J=600; K=700; L=750; M=725; O="Provincia de Diguillín"; P=725},
    @{Row=114; D=44376; I="Primera"; J=120; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=115; D=44172; I="Primera"; J=160; K=700; L=750; M=725; O="Provincia de Diguillín"; P=725},
    @{Row=116; D=44421; I="Primera"; J=300; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=117; D=44426; I="Primera"; J=300; K=500; L=550; M=525; O="Provincia de Diguillín"; P=525},
    @{Row=118; D=44448; I="Primera"; J=300; K=600; L=650; M=625; O="Región del Maule"; P=625},
    @{Row=119; D=44362; I="Primera"; J=120; K=500; L=550; M=525; O="Provincia de Diguillín"; P=525},
    @{Row=120; D=44210; I="Primera"; J=2600; K=700; L=800; M=750; O="Provincia de Diguillín"; P=750},
    @{Row=121; D=44176; I="Primera"; J=300; K=800; L=850; M=825; O="Provincia de Diguillín"; P=825},
    @{Row=122; D=44301; I="Primera"; J=300; K=800; L=850; M=825; O="Provincia de Diguillín"; P=825},
    @{Row=123; D=44407; I="Primera"; J=300; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=124; D=44284; I="Primera"; J=120; K=700; L=750; M=725; O="Provincia de Diguillín"; P=725},
    @{Row=125; D=44441; I="Primera"; J=300; K=600; L=650; M=625; O="Región del Maule"; P=625},
    @{Row=126; D=44504; I="Primera"; J=360; K=600; L=700; M=650; O="Provincia de Diguillín"; P=650},
    @{Row=127; D=44350; I="Primera"; J=300; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=128; D=44382; I="Primera"; J=160; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=129; D=44329; I="Primera"; J=300; K=500; L=550; M=525; O="Provincia de Diguillín"; P=525},
    @{Row=130; D=44491; I="Primera"; J=200; K=600; L=700; M=650; O="Región del Maule"; P=650},
    @{Row=131; D=44305; I="Primera"; J=120; K=800; L=850; M=825; O="Provincia de Diguillín"; P=825},
    @{Row=132; D=44225; I="Primera"; J=2800; K=900; L=1000; M=946; O="Región del Maule"; P=946},
    @{Row=133; D=44447; I="Primera"; J=300; K=700; L=750; M=725; O="Provincia de Diguillín"; P=725},
    @{Row=134; D=44425; I="Primera"; J=300; K=500; L=550; M=525; O="Provincia de Diguillín"; P=525},
    @{Row=135; D=44322; I="Primera"; J=600; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=136; D=44495; I="Primera"; J=180; K=600; L=700; M=650; O="Región del Maule"; P=650},
    @{Row=137; D=44232; I="Primera"; J=300; K=800; L=850; M=825; O="Provincia de Diguillín"; P=825},
    @{Row=138; D=44327; I="Primera"; J=300; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=139; D=44510; I="Primera"; J=360; K=600; L=700; M=650; O="Provincia de Diguillín"; P=650},
    @{Row=140; D=44161; I="Primera"; J=2600; K=950; L=1000; M=977; O="Región del Maule"; P=977},
    @{Row=141; D=44468; I="Primera"; J=300; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=142; D=44517; I="Primera"; J=200; K=600; L=700; M=650; O="Provincia de Diguillín"; P=650},
    @{Row=143; D=44391; I="Primera"; J=160; K=500; L=600; M=550; O="Provincia de Diguillín"; P=550},
    @{Row=144; D=44236; I="Primera"; J=300; K=700; L=750; M=725; O="Provincia de Diguillín"; P=725},
    @{Row=145; D=44340; I="Primera"; J=160; K=700; L=750; M=725; O="Provincia de Diguillín"; P=725},
    @{Row=146; D=44515; I="Primera"; J=240; K=600; L=700; M=650; O="Provincia de Diguillín"; P=650},
    @{Row=147; D=44330; I="Primera"; J=300; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=148; D=44432; I="Primera"; J=300; K=600; L=650; M=625; O="Provincia de Diguillín"; P=625},
    @{Row=149; D=44181; I="Primera"; J=160; K=800; L=850; M=825; O="Provincia de Diguillín"; P=825},
    @{Row=150; D=44194; I="Primera"; J=160; K=700; L=750; M=725; O="Provincia de Diguillín"; P=725},
    @{Row=151; D=44307; I="Primera"; J=160; K=700; L=750; M=725; O="Provincia de Diguillín"; P=725},
    @{Row=152; D=44508; I="Primera"; J=400; K=600; L=700; M=650; O="Provincia de Diguillín"; P=650}
)

$colIndex = @{ D = 4; I = 9; J = 10; K = 11; L = 12; M = 13; O = 15; P = 16 }

foreach ($rec in $records) {
    $r = $rec.Row
    foreach ($col in $colIndex.Keys) {
        if ($rec.ContainsKey($col)) {
            $ws.Cells.Item($r, $colIndex[$col]).Value = $rec[$col]
        }
    }
}

# New row 152 also needs the columns that stay constant across the whole sheet
# (Mercado ID, Mercado, Region, Codreg, Categoria ID, Categoria, Variedad,
# Unidad de comercializacion, Kg o Unidades, Clasificacion), copied from row 151.
$ws.Cells.Item(152, 1).Value = 7
$ws.Cells.Item(152, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(152, 3).Value = "Ñuble"
$ws.Cells.Item(152, 5).Value = 16
$ws.Cells.Item(152, 6).Value = 100112006
$ws.Cells.Item(152, 7).Value = "Repollo"
$ws.Cells.Item(152, 8).Value = "Crespo record"
$ws.Cells.Item(152, 14).Value = "`$/unidad"
$ws.Cells.Item(152, 17).Value = 1
$ws.Cells.Item(152, 18).Value = "Hortaliza"

# Row 152, column D is a date cell; copy the date number format from row 151
# so the new cell gets the same style index as the rest of the Fecha column.
$ws.Cells.Item(152, 4).NumberFormat = $ws.Cells.Item(151, 4).NumberFormat

Write-Host "Applied weekly shift to rows 47-152."
